$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1581976666666667
$ws.Range("H2").Value = 0.474593
$ws.Range("I2").Value = 0.1400666049254827
$ws.Range("J2").Value = 0.1400666049254826
$ws.Range("M2").Value = 0.166446
$ws.Range("N2").Value = 0.4993379999999999
$ws.Range("O2").Value = 0.008607068890887148
$ws.Range("P2").Value = 0.008607068890887149
$ws.Range("Q2").Value = 0.026331368826
$ws.Range("R2").Value = 0.236982319434
$ws.Range("S2").Value = 0.001205562917906302
$ws.Range("T2").Value = 0.001205562917906302
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1581976666666667
$ws.Range("H3").Value = 0.474593
$ws.Range("I3").Value = 0.1400666049254827
$ws.Range("J3").Value = 0.1400666049254826
$ws.Range("O3").Value = 0.09774836700492934
$ws.Range("P3").Value = 0.09774836700492935
$ws.Range("Q3").Value = 0.2990388872652223
$ws.Range("R3").Value = 2.691349985387
$ws.Range("S3").Value = 0.01369128190339052
$ws.Range("T3").Value = 0.01369128190339052
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1581976666666667
$ws.Range("H4").Value = 0.474593
$ws.Range("I4").Value = 0.1400666049254827
$ws.Range("J4").Value = 0.1400666049254826
$ws.Range("M4").Value = 17.13435266666667
$ws.Range("N4").Value = 51.403058
$ws.Range("O4").Value = 0.8860324297535294
$ws.Range("P4").Value = 0.8860324297535294
$ws.Range("Q4").Value = 2.710614611710445
$ws.Range("R4").Value = 24.395531505394
$ws.Range("S4").Value = 0.1241035542894531
$ws.Range("T4").Value = 0.124103554289453
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1581976666666667
$ws.Range("H5").Value = 0.474593
$ws.Range("I5").Value = 0.1400666049254827
$ws.Range("J5").Value = 0.1400666049254826
$ws.Range("M5").Value = 0.1472056666666667
$ws.Range("N5").Value = 0.441617
$ws.Range("O5").Value = 0.007612134350654087
$ws.Range("P5").Value = 0.007612134350654088
$ws.Range("Q5").Value = 0.02328759298677778
$ws.Range("R5").Value = 0.209588336881
$ws.Range("S5").Value = 0.001066205814732761
$ws.Range("T5").Value = 0.001066205814732761
$ws.Range("I6").Value = 0.7029419733214338
$ws.Range("J6").Value = 0.7029419733214337
$ws.Range("M6").Value = 0.166446
$ws.Range("N6").Value = 0.4993379999999999
$ws.Range("O6").Value = 0.008607068890887148
$ws.Range("P6").Value = 0.008607068890887149
$ws.Range("Q6").Value = 0.13214730501
$ws.Range("R6").Value = 1.18932574509
$ws.Range("S6").Value = 0.006050269990673736
$ws.Range("T6").Value = 0.006050269990673736
$ws.Range("I7").Value = 0.7029419733214338
$ws.Range("J7").Value = 0.7029419733214337
$ws.Range("O7").Value = 0.09774836700492934
$ws.Range("P7").Value = 0.09774836700492935
$ws.Range("S7").Value = 0.06871142999139276
$ws.Range("T7").Value = 0.06871142999139276
$ws.Range("I8").Value = 0.7029419733214338
$ws.Range("J8").Value = 0.7029419733214337
$ws.Range("M8").Value = 17.13435266666667
$ws.Range("N8").Value = 51.403058
$ws.Range("O8").Value = 0.8860324297535294
$ws.Range("P8").Value = 0.8860324297535294
$ws.Range("Q8").Value = 13.60356228441
$ws.Range("R8").Value = 122.43206055969
$ws.Range("S8").Value = 0.6228293845977306
$ws.Range("T8").Value = 0.6228293845977305
$ws.Range("I9").Value = 0.7029419733214338
$ws.Range("J9").Value = 0.7029419733214337
$ws.Range("M9").Value = 0.1472056666666667
$ws.Range("N9").Value = 0.441617
$ws.Range("O9").Value = 0.007612134350654087
$ws.Range("P9").Value = 0.007612134350654088
$ws.Range("Q9").Value = 0.116871730965
$ws.Range("R9").Value = 1.051845578685
$ws.Range("S9").Value = 0.005350888741636656
$ws.Range("T9").Value = 0.005350888741636655
$ws.Range("G10").Value = 0.1199896666666667
$ws.Range("H10").Value = 0.359969
$ws.Range("I10").Value = 0.106237630366274
$ws.Range("J10").Value = 0.106237630366274
$ws.Range("M10").Value = 0.166446
$ws.Range("N10").Value = 0.4993379999999999
$ws.Range("O10").Value = 0.008607068890887148
$ws.Range("P10").Value = 0.008607068890887149
$ws.Range("Q10").Value = 0.019971800058
$ws.Range("R10").Value = 0.179746200522
$ws.Range("S10").Value = 0.0009143946033671246
$ws.Range("T10").Value = 0.0009143946033671245
$ws.Range("G11").Value = 0.1199896666666667
$ws.Range("H11").Value = 0.359969
$ws.Range("I11").Value = 0.106237630366274
$ws.Range("J11").Value = 0.106237630366274
$ws.Range("O11").Value = 0.09774836700492934
$ws.Range("P11").Value = 0.09774836700492935
$ws.Range("Q11").Value = 0.2268148270412222
$ws.Range("R11").Value = 2.041333443371
$ws.Range("S11").Value = 0.01038455488277657
$ws.Range("T11").Value = 0.01038455488277657
$ws.Range("G12").Value = 0.1199896666666667
$ws.Range("H12").Value = 0.359969
$ws.Range("I12").Value = 0.106237630366274
$ws.Range("J12").Value = 0.106237630366274
$ws.Range("M12").Value = 17.13435266666667
$ws.Range("N12").Value = 51.403058
$ws.Range("O12").Value = 0.8860324297535294
$ws.Range("P12").Value = 0.8860324297535294
$ws.Range("Q12").Value = 2.055945265022444
$ws.Range("R12").Value = 18.503507385202
$ws.Range("S12").Value = 0.09412998576468708
$ws.Range("T12").Value = 0.09412998576468705
$ws.Range("G13").Value = 0.1199896666666667
$ws.Range("H13").Value = 0.359969
$ws.Range("I13").Value = 0.106237630366274
$ws.Range("J13").Value = 0.106237630366274
$ws.Range("M13").Value = 0.1472056666666667
$ws.Range("N13").Value = 0.441617
$ws.Range("O13").Value = 0.007612134350654087
$ws.Range("P13").Value = 0.007612134350654088
$ws.Range("Q13").Value = 0.01766315887477778
$ws.Range("R13").Value = 0.158968429873
$ws.Range("S13").Value = 0.0008086951154432059
$ws.Range("T13").Value = 0.0008086951154432058
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05732366666666666
$ws.Range("H14").Value = 0.171971
$ws.Range("I14").Value = 0.05075379138680971
$ws.Range("J14").Value = 0.05075379138680969
$ws.Range("M14").Value = 0.166446
$ws.Range("N14").Value = 0.4993379999999999
$ws.Range("O14").Value = 0.008607068890887148
$ws.Range("P14").Value = 0.008607068890887149
$ws.Range("Q14").Value = 0.009541295021999999
$ws.Range("R14").Value = 0.08587165519799998
$ws.Range("S14").Value = 0.0004368413789399859
$ws.Range("T14").Value = 0.0004368413789399859
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05732366666666666
$ws.Range("H15").Value = 0.171971
$ws.Range("I15").Value = 0.05075379138680971
$ws.Range("J15").Value = 0.05075379138680969
$ws.Range("O15").Value = 0.09774836700492934
$ws.Range("P15").Value = 0.09774836700492935
$ws.Range("Q15").Value = 0.1083581436765555
$ws.Range("R15").Value = 0.9752232930889999
$ws.Range("S15").Value = 0.004961100227369497
$ws.Range("T15").Value = 0.004961100227369496
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05732366666666666
$ws.Range("H16").Value = 0.171971
$ws.Range("I16").Value = 0.05075379138680971
$ws.Range("J16").Value = 0.05075379138680969
$ws.Range("M16").Value = 17.13435266666667
$ws.Range("N16").Value = 51.403058
$ws.Range("O16").Value = 0.8860324297535294
$ws.Range("P16").Value = 0.8860324297535294
$ws.Range("Q16").Value = 0.9822039208131111
$ws.Range("R16").Value = 8.839835287318
$ws.Range("S16").Value = 0.04496950510165876
$ws.Range("T16").Value = 0.04496950510165874
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05732366666666666
$ws.Range("H17").Value = 0.171971
$ws.Range("I17").Value = 0.05075379138680971
$ws.Range("J17").Value = 0.05075379138680969
$ws.Range("M17").Value = 0.1472056666666667
$ws.Range("N17").Value = 0.441617
$ws.Range("O17").Value = 0.007612134350654087
$ws.Range("P17").Value = 0.007612134350654088
$ws.Range("Q17").Value = 0.008438368567444444
$ws.Range("R17").Value = 0.075945317107
$ws.Range("S17").Value = 0.0003863446788414657
$ws.Range("T17").Value = 0.0003863446788414656
